# Protons-vs-ions.xlsx : add a "Carbon" ion reference-particle block next to the
# existing "Proton" one, wire the ion-parameters block to it, fix a couple of
# labels, and tweak the drift length used for the proton reference particle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Proton reference-particle block: drift length 20 -> 15 MeV
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 15

# ---------------------------------------------------------------------------
# 2. Build the new "Carbon" block in columns L:N (mirrors A:C for the proton
#    block) by copying the formatting across first, then filling in content.
# ---------------------------------------------------------------------------
$ws.Range("A2:C8").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("L2:N2").Merge()

$ws.Range("L1").Value = "Carbon"

$ws.Range("L2").Value = "Reference particle"

$ws.Range("L3").Value = "Kr"
$ws.Range("M3").Formula = "=M4-I11"
$ws.Range("N3").Value = "MeV"

$ws.Range("L4").Value = "Er"
$ws.Range("M4").Formula = "=SQRT(I11^2+M5^2)"
$ws.Range("N4").Value = "MeV"

$ws.Range("L5").Value = "pr"
$ws.Range("M5").Formula = "=B14*I12"
$ws.Range("N5").Value = "MeV"

$ws.Range("L6").Value = "grbr"
$ws.Range("M6").Formula = "=M5/T3"

$ws.Range("L7").Value = "br"
$ws.Range("M7").Formula = "=M5/M4"

$ws.Range("L8").Value = "charge"
$ws.Range("M8").Formula = "=I12"

$ws.Range("L9").Value = "Kr"
$ws.Range("M9").Formula = "=M3/12"
$ws.Range("N9").Value = "MeV/u"

# ---------------------------------------------------------------------------
# 3. Ion-parameters block (H10:J12 / H14:J16 / H22:I22): charge becomes a
#    formula (6 * |e|-unit), and the kinetic-energy figure is now sourced
#    from the new Carbon block instead of being computed in place.
# ---------------------------------------------------------------------------
$ws.Range("I12").Formula = "=6*I2"
$ws.Range("J12").ClearContents()

$ws.Range("H9").Value = "Ion parameters"
$ws.Range("H14").Value = "KE"
$ws.Range("I14").Formula = "=M3"

$ws.Range("I22").Formula = "=I16/I12"

# ---------------------------------------------------------------------------
# 4. Small standalone MeV -> MeV/u conversion check near the bottom.
# ---------------------------------------------------------------------------
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = "MeV/u"
$ws.Range("F27").Formula = "=F26*12"

# ---------------------------------------------------------------------------
# 5. Restore selection close to where the author left it.
# ---------------------------------------------------------------------------
$ws.Range("L17:M20").Select()
